$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 224.63637
$ws.Range("J2").Value = 326.66666
$ws.Range("L2").Value = 326.66666
$ws.Range("N2").Value = -552.66666
$ws.Range("H15").Value = 94150.45
$ws.Range("I15").Value = 94150.45
$ws.Range("K15").Value = 282451.35
$ws.Range("M15").Value = -282282.35
$ws.Range("H53").Value = 297.7
$ws.Range("I53").Value = 161.07692
$ws.Range("J53").Value = 551.4286
$ws.Range("K53").Value = 161.07692
$ws.Range("L53").Value = 551.4286
$ws.Range("M53").Value = 475.92308
$ws.Range("N53").Value = -1825.4286
$ws.Range("H62").Value = 2941673.8
$ws.Range("I62").Value = 4129058.5
$ws.Range("J62").Value = 27184.273
$ws.Range("K62").Value = 4129058.5
$ws.Range("L62").Value = 27184.273
$ws.Range("M62").Value = -4128434.5
$ws.Range("N62").Value = -28432.273
$ws.Range("H65").Value = 2941673.8
$ws.Range("I65").Value = 4129058.5
$ws.Range("J65").Value = 27184.273
$ws.Range("K65").Value = 20645292.5
$ws.Range("L65").Value = 135921.365
$ws.Range("M65").Value = -20642172.5
$ws.Range("N65").Value = -142161.365
$ws.Range("H76").Value = 3370063.5
$ws.Range("I76").Value = 3971146.2
$ws.Range("K76").Value = 3971146.2
$ws.Range("M76").Value = -3970831.2
$ws.Range("H79").Value = 3370063.5
$ws.Range("I79").Value = 3971146.2
$ws.Range("K79").Value = 3971146.2
$ws.Range("M79").Value = -3970054.2
$ws.Range("H111").Value = 1271.9375
$ws.Range("I111").Value = 718.4286
$ws.Range("J111").Value = 1702.4445
$ws.Range("K111").Value = 2155.2858
$ws.Range("L111").Value = 5107.333500000001
$ws.Range("M111").Value = 911.7142000000003
$ws.Range("N111").Value = -11241.3335
$ws.Range("H132").Value = 30341.371
$ws.Range("I132").Value = 34968.965
$ws.Range("J132").Value = 2575.8
$ws.Range("K132").Value = 104906.895
$ws.Range("L132").Value = 7727.400000000001
$ws.Range("M132").Value = -102376.895
$ws.Range("N132").Value = -12787.4
$ws.Range("H133").Value = 45590.527
$ws.Range("J133").Value = 45590.527
$ws.Range("L133").Value = 45590.527
$ws.Range("N133").Value = -55710.527
$ws.Range("H137").Value = 23256738
$ws.Range("I137").Value = 26316504
$ws.Range("J137").Value = 2500.6
$ws.Range("K137").Value = 78949512
$ws.Range("L137").Value = 7501.799999999999
$ws.Range("M137").Value = -78946962
$ws.Range("N137").Value = -12601.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 3650
$ws.Range("I31").Value = 3650
$ws.Range("K31").Value = 3650
$ws.Range("M31").Value = -3356
$ws.Range("H32").Value = 18651.676
$ws.Range("I32").Value = 3673.2239
$ws.Range("J32").Value = 95848.30499999999
$ws.Range("K32").Value = 3673.2239
$ws.Range("L32").Value = 95848.30499999999
$ws.Range("M32").Value = -3386.2239
$ws.Range("N32").Value = -96422.30499999999
$ws.Range("H45").Value = 813.38464
$ws.Range("I45").Value = 813.38464
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 813.38464
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -436.38464
$ws.Range("N45").Value = $null
$ws.Range("H61").Value = 2991.4546
$ws.Range("I61").Value = 1594.909
$ws.Range("J61").Value = 4388
$ws.Range("K61").Value = 1594.909
$ws.Range("L61").Value = 4388
$ws.Range("M61").Value = -1382.909
$ws.Range("N61").Value = -4812
$ws.Range("H74").Value = 4976.25
$ws.Range("I74").Value = 1159
$ws.Range("K74").Value = 1159
$ws.Range("M74").Value = -285
$ws.Range("H77").Value = 4976.25
$ws.Range("I77").Value = 1159
$ws.Range("K77").Value = 5795
$ws.Range("M77").Value = -1427
$ws.Range("H122").Value = 12281.1
$ws.Range("I122").Value = 13423.444
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 40270.33199999999
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -37820.33199999999
$ws.Range("N122").Value = -10900
$ws.Range("H136").Value = 2991.4546
$ws.Range("I136").Value = 1594.909
$ws.Range("J136").Value = 4388
$ws.Range("K136").Value = 4784.727000000001
$ws.Range("L136").Value = 13164
$ws.Range("M136").Value = -2234.727000000001
$ws.Range("N136").Value = -18264
$ws.Range("H139").Value = 53715
$ws.Range("J139").Value = 53715
$ws.Range("L139").Value = 53715
$ws.Range("N139").Value = -63995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = $null
$ws.Range("H134").Value = 3206.375
$ws.Range("I134").Value = 2419.4443
$ws.Range("J134").Value = 4840.769
$ws.Range("K134").Value = 7258.3329
$ws.Range("L134").Value = 14522.307
$ws.Range("M134").Value = -4723.3329
$ws.Range("N134").Value = -19592.307
$ws.Range("H139").Value = 85390
$ws.Range("J139").Value = 85390
$ws.Range("L139").Value = 85390
$ws.Range("N139").Value = -95670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1827.9546
$ws.Range("I58").Value = 1268.2354
$ws.Range("K58").Value = 1268.2354
$ws.Range("M58").Value = -1065.2354
$ws.Range("H70").Value = 29000
$ws.Range("J70").Value = 29000
$ws.Range("L70").Value = 29000
$ws.Range("N70").Value = -29630
$ws.Range("H73").Value = 29000
$ws.Range("J73").Value = 29000
$ws.Range("L73").Value = 29000
$ws.Range("N73").Value = -31184
$ws.Range("H132").Value = 1579.6364
$ws.Range("I132").Value = 964.4231
$ws.Range("K132").Value = 2893.2693
$ws.Range("M132").Value = -363.2692999999999
$ws.Range("H134").Value = 2392.2563
$ws.Range("I134").Value = 1457.5454
$ws.Range("J134").Value = 7533.1665
$ws.Range("K134").Value = 4372.6362
$ws.Range("L134").Value = 22599.4995
$ws.Range("M134").Value = -1837.6362
$ws.Range("N134").Value = -27669.4995
$ws.Range("H136").Value = 1827.9546
$ws.Range("I136").Value = 1268.2354
$ws.Range("K136").Value = 3804.7062
$ws.Range("M136").Value = -1254.7062

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2300.8462
$ws.Range("I5").Value = 1647.091
$ws.Range("J5").Value = 2557.6785
$ws.Range("K5").Value = 4941.272999999999
$ws.Range("L5").Value = 7673.0355
$ws.Range("M5").Value = -4829.272999999999
$ws.Range("N5").Value = -7897.0355
$ws.Range("H107").Value = 465.92307
$ws.Range("I107").Value = 455.88235
$ws.Range("J107").Value = 484.8889
$ws.Range("K107").Value = 1367.64705
$ws.Range("L107").Value = 1454.6667
$ws.Range("M107").Value = 552.35295
$ws.Range("N107").Value = -5294.6667
$ws.Range("H131").Value = 4167936
$ws.Range("J131").Value = 4696217
$ws.Range("L131").Value = 14088651
$ws.Range("N131").Value = -14098731
$ws.Range("H135").Value = 2300.8462
$ws.Range("I135").Value = 1647.091
$ws.Range("J135").Value = 2557.6785
$ws.Range("K135").Value = 14823.819
$ws.Range("L135").Value = 23019.1065
$ws.Range("M135").Value = -12288.819
$ws.Range("N135").Value = -28089.1065

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 24400
$ws.Range("J95").Value = 24400
$ws.Range("L95").Value = 24400
$ws.Range("N95").Value = -29892
$ws.Range("H102").Value = 2272.0833
$ws.Range("I102").Value = 1892.7142
$ws.Range("K102").Value = 1892.7142
$ws.Range("M102").Value = -270.7141999999999
$ws.Range("H137").Value = 50226.332
$ws.Range("J137").Value = 50226.332
$ws.Range("L137").Value = 50226.332
$ws.Range("N137").Value = -60426.332
$ws.Range("H138").Value = 76809.664
$ws.Range("J138").Value = 76809.664
$ws.Range("L138").Value = 76809.664
$ws.Range("N138").Value = -87089.664
$ws.Range("H139").Value = 43870
$ws.Range("J139").Value = 43870
$ws.Range("L139").Value = 43870
$ws.Range("N139").Value = -54150

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2594.8206
$ws.Range("I7").Value = 1935.8235
$ws.Range("J7").Value = 3104.0454
$ws.Range("K7").Value = 1935.8235
$ws.Range("L7").Value = 3104.0454
$ws.Range("M7").Value = -1823.8235
$ws.Range("N7").Value = -3328.0454
$ws.Range("H40").Value = 3182.5217
$ws.Range("I40").Value = 2116.3333
$ws.Range("J40").Value = 3558.8235
$ws.Range("K40").Value = 2116.3333
$ws.Range("L40").Value = 3558.8235
$ws.Range("M40").Value = -1980.3333
$ws.Range("N40").Value = -3830.8235
$ws.Range("H93").Value = 1452
$ws.Range("I93").Value = 1239.6666
$ws.Range("J93").Value = 4000
$ws.Range("K93").Value = 1239.6666
$ws.Range("L93").Value = 4000
$ws.Range("M93").Value = 8.333399999999983
$ws.Range("N93").Value = -6496
$ws.Range("H101").Value = 29999.666
$ws.Range("J101").Value = 29999.666
$ws.Range("L101").Value = 29999.666
$ws.Range("N101").Value = -36489.666
$ws.Range("H126").Value = 2594.8206
$ws.Range("I126").Value = 1935.8235
$ws.Range("J126").Value = 3104.0454
$ws.Range("K126").Value = 5807.470499999999
$ws.Range("L126").Value = 9312.136200000001
$ws.Range("M126").Value = -3337.470499999999
$ws.Range("N126").Value = -14252.1362
$ws.Range("H132").Value = 3772.1191
$ws.Range("I132").Value = 2225.5925
$ws.Range("J132").Value = 6555.8667
$ws.Range("K132").Value = 6676.7775
$ws.Range("L132").Value = 19667.6001
$ws.Range("M132").Value = -4146.7775
$ws.Range("N132").Value = -24727.6001
$ws.Range("H136").Value = 4488.1816
$ws.Range("I136").Value = 1881.9445
$ws.Range("J136").Value = 16216.25
$ws.Range("K136").Value = 5645.833500000001
$ws.Range("L136").Value = 48648.75
$ws.Range("M136").Value = -3095.833500000001
$ws.Range("N136").Value = -53748.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1398.7273
$ws.Range("I122").Value = 1585.875
$ws.Range("J122").Value = 899.6667
$ws.Range("K122").Value = 4757.625
$ws.Range("L122").Value = 2699.0001
$ws.Range("M122").Value = -2307.625
$ws.Range("N122").Value = -7599.0001
$ws.Range("H126").Value = 54939.473
$ws.Range("I126").Value = 85795.836
$ws.Range("J126").Value = 2042.8572
$ws.Range("K126").Value = 257387.508
$ws.Range("L126").Value = 6128.571599999999
$ws.Range("M126").Value = -254917.508
$ws.Range("N126").Value = -11068.5716
$ws.Range("H136").Value = 3266.2703
$ws.Range("I136").Value = 1187.5385
$ws.Range("J136").Value = 8179.636
$ws.Range("K136").Value = 3562.6155
$ws.Range("L136").Value = 24538.908
$ws.Range("M136").Value = -1012.6155
$ws.Range("N136").Value = -29638.908
